$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (Albatros / caballeros entry), pushing
# existing rows 2-6 down to 3-7.
$ws.Rows.Item(2).Insert()

# Insert a new row at what is now position 5 (Prejuveniles / caballeros
# second entry), pushing the remaining rows down further.
$ws.Rows.Item(5).Insert()

# Row 2 - new entry: Albatros / caballeros
$torneo = "####1er Torneo Federativo - C.A.E. - Sub 23, Prejuveniles y sub 23 (28 de Febrero y 1 de Marzo) - Juniors (Domingo 1 de Marzo)"
$ws.Cells.Item(2, 1).Value = $torneo
$ws.Cells.Item(2, 2).Value = "Albatros"
$ws.Cells.Item(2, 3).Value = "caballeros"
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = "Luján Martínez, Benjamín"
$ws.Cells.Item(2, 6).Value = 98
$ws.Cells.Item(2, 7).Value = ""
$ws.Cells.Item(2, 8).Value = 98

# Row 5 - new entry: Prejuveniles / caballeros (position 2), value in dia_2
$ws.Cells.Item(5, 1).Value = $torneo
$ws.Cells.Item(5, 2).Value = "Prejuveniles"
$ws.Cells.Item(5, 3).Value = "caballeros"
$ws.Cells.Item(5, 4).Value = 2
$ws.Cells.Item(5, 5).Value = "Luján Martínez, Benjamín"
$ws.Cells.Item(5, 6).Value = ""
$ws.Cells.Item(5, 7).Value = 98
$ws.Cells.Item(5, 8).Value = 98

$wb.Save()
